# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.633.26'
$ws.Range('E2').Value = '  -3.09%  '

$ws.Range('D3').Value = '2.982.04'
$ws.Range('E3').Value = '  -5.57%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '496.83'
$ws.Range('E5').Value = '  -5.95%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.04'
$ws.Range('E6').Value = '  +0.24%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('E8').Value = '  -4.54%  '

$ws.Range('E9').Value = '  -1.28%  '

$ws.Range('E10').Value = '  -4.00%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.350'
$ws.Range('E11').Value = '  -7.50%  '

$ws.Range('E12').Value = '  -0.75%  '

$ws.Range('D13').Value = '3.492.96'
$ws.Range('E13').Value = '  -5.48%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.26'
$ws.Range('E14').Value = '  -0.75%  '

$ws.Range('D15').Value = '56.643.77'
$ws.Range('E15').Value = '  -2.97%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000147'
$ws.Range('E16').Value = '  -3.30%  '

$ws.Range('D17').Value = '2.982.50'
$ws.Range('E17').Value = '  -5.51%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.71'
$ws.Range('E18').Value = '  -1.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.35'
$ws.Range('E19').Value = '  -5.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.77'
$ws.Range('E20').Value = '  -2.55%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.91'
$ws.Range('E21').Value = '  -5.39%  '

$ws.Range('E22').Value = '  -0.08%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.469'
$ws.Range('E23').Value = '  -8.26%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.67'
$ws.Range('E24').Value = '  -8.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.35%  '

$ws.Range('E26').Value = '  -5.94%  '

$ws.Range('D27').Value = '0.0₃0893'
$ws.Range('E27').Value = '  -6.37%  '

$ws.Range('E28').Value = '  -0.05%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.47'
$ws.Range('E29').Value = '  -5.57%  '

$ws.Range('E30').Value = '  -2.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.74'
$ws.Range('E31').Value = '  -7.27%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.22'
$ws.Range('E32').Value = '  -5.98%  '

$ws.Range('E33').Value = '  -7.92%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.64'
$ws.Range('E34').Value = '  -4.36%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.43'
$ws.Range('E35').Value = '  -8.70%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.27'
$ws.Range('E36').Value = '  -7.92%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.59'
$ws.Range('E37').Value = '  -10.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0669'
$ws.Range('E38').Value = '  -2.73%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.33'
$ws.Range('E39').Value = '  -2.73%  '

$ws.Range('D40').Value = '3.016.07'
$ws.Range('E40').Value = '  -5.37%  '

$ws.Range('E41').Value = '  -9.40%  '

$ws.Range('E42').Value = '  +0.20%  '

$ws.Range('E43').Value = '  -8.21%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.997'
$ws.Range('E44').Value = '  -8.53%  '

$ws.Range('D45').Value = '2.221.41'
$ws.Range('E45').Value = '  -3.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.40'
$ws.Range('E46').Value = '  -4.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.55'
$ws.Range('E47').Value = '  -10.02%  '

$ws.Range('E48').Value = '  +4.70%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0235'
$ws.Range('E49').Value = '  +0.19%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.72'
$ws.Range('E50').Value = '  -7.46%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.94'
$ws.Range('E51').Value = '  -8.59%  '
